$wb = $excel.ActiveWorkbook

# The "SOmixos" sheet (second sheet) has a "Type" column (D) whose values
# ("PreyCurve"/"Treatments") are redundant with other columns. Remove the
# whole column, shifting everything after it one column to the left.
$ws = $wb.Worksheets.Item("SOmixos")
$ws.Columns("D").Delete()

# Select column D (now GrazingRate) on SOmixos and make it the active sheet/tab.
$ws.Select() | Out-Null
$ws.Range("D1:D1048576").Select() | Out-Null
